$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 83.33333333333334
$ws.Range("C2").Value = 83.33333333333334
$ws.Range("D2").Value = 81.77083333333334
$ws.Range("E2").Value = 59.89583333333333
$ws.Range("F2").Value = 53.125
$ws.Range("G2").Value = 52.86458333333333
$ws.Range("H2").Value = 51.82291666666667
$ws.Range("I2").Value = 51.30208333333333
